$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C15: achievements_rate_per_100000_population caveat text - clarify naming of per 100k rate
$ws.Range("C15").Value = "The rates are the number of achievments in AY21/22 per 100,000 of the population. Further education and skills include all age apprenticeships and publicly-funded adult (19+) learning, including community learning, delivered by an FE institution, a training provider or within a local community. `n"

# C16: participation_rate_per_100000_population caveat text - clarify naming of per 100k rate
$ws.Range("C16").Value = "The rates are the number of participants in AY21/22 per 100,000 of the population. Further education and skills include all age apprenticeships and publicly-funded adult (19+) learning, including community learning, delivered by an FE institution, a training provider or within a local community. `n"

# Leave the cursor where the author ended up after editing these rows
$ws.Range("A17").Select() | Out-Null
